# Append new chat-log rows (51-69) to the log sheet: UserID, UserInput,
# BotReply, Timestamp. Rows 53/55 intentionally have no BotReply (bus
# location lookups that returned nothing), matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("U53e8616065a61b190a510f520198a271", "bus KKA-8155", "車輛:KKA-8155`n業者代號:34`nGPS車速:98.0`nGPS時間:2022-05-09 11:18:59`n路線方向(1:去程,2:回程):1", "2022-05-09 22:19:29.557000")
    ,@("U53e8616065a61b190a510f520198a271", "bus KKA-8155", "車輛:KKA-8155`n業者代號:34`nGPS車速:98.0`nGPS時間:2022-05-09 11:18:59`n路線方向(1:去程,2:回程):1", "2022-05-09 22:22:56.174000")
    ,@("U53e8616065a61b190a510f520198a271", "bus fuck", "", "2022-05-09 22:25:34.543000")
    ,@("U53e8616065a61b190a510f520198a271", "bus KKA-8155", "車輛:KKA-8155`n業者代號:34`nGPS車速:98.0`nGPS時間:2022-05-09 11:18:59`n路線方向(1:去程,2:回程):1", "2022-05-09 22:41:54.154000")
    ,@("U53e8616065a61b190a510f520198a271", "bus. 123", "", "2022-05-09 22:42:35.137000")
    ,@("U53e8616065a61b190a510f520198a271", "bus 123", "桃園公車中查無此資料", "2022-05-10 09:29:19.487000")
    ,@("U53e8616065a61b190a510f520198a271", "筆記本", "筆記本50元", "2022-05-10 09:29:29.618000")
    ,@("U53e8616065a61b190a510f520198a271", "地址", "261宜蘭縣頭城鎮港口路92-1號", "2022-05-10 09:29:53.618000")
    ,@("U53e8616065a61b190a510f520198a271", "ubike 健行科技大學", "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:15`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1", "2022-05-10 09:30:08.741000")
    ,@("U53e8616065a61b190a510f520198a271", "ubike 健行科技大學", "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:15`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1", "2022-05-10 09:53:14.429000")
    ,@("U53e8616065a61b190a510f520198a271", "ubike 健行科技大學", "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:15`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1", "2022-05-10 09:55:03.262000")
    ,@("U53e8616065a61b190a510f520198a271", "地址", "261宜蘭縣頭城鎮港口路92-1號", "2022-05-10 09:55:29.282000")
    ,@("U53e8616065a61b190a510f520198a271", "bus KKA-8155", "車輛:KKA-8155`n業者代號:34`nGPS車速:98.0`nGPS時間:2022-05-09 11:18:59`n路線方向(1:去程,2:回程):1", "2022-05-10 09:55:43.636000")
    ,@("U53e8616065a61b190a510f520198a271", "我", "Unable to recognize user's input", "2022-05-10 09:58:55.522000")
    ,@("U53e8616065a61b190a510f520198a271", "地址", "261宜蘭縣頭城鎮港口路92-1號", "2022-05-10 10:48:39.557000")
    ,@("U53e8616065a61b190a510f520198a271", "ubike 健行科技大學", "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:15`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1", "2022-05-10 10:48:51.925000")
    ,@("U53e8616065a61b190a510f520198a271", "bus KKA-8155", "車輛:KKA-8155`n業者代號:34`nGPS車速:98.0`nGPS時間:2022-05-09 11:18:59`n路線方向(1:去程,2:回程):1", "2022-05-10 10:49:05.773000")
    ,@("U53e8616065a61b190a510f520198a271", "有表單嗎", "https://forms.gle/vdHfmWijtcBTsPNX6", "2022-05-10 10:55:14.876000")
    ,@("U53e8616065a61b190a510f520198a271", "位置情報", "nontextreply", "2022-05-10 11:23:50.980000")
)

$startRow = 51
$endRow = $startRow + $data.Count - 1
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    if ($vals[2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $vals[2]
    }
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$ws.Range("A" + $startRow + ":A" + $endRow).EntireRow.AutoFit()
